$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2:H9").Value = 2.5
$ws.Range("H4:H9").Borders.LineStyle = -4142

$ws.Range("H9").Select()
